$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "5,53"
$ws.Range("F2").Value = "11 de jun., 16:37 UTC ·"

# Row 3
$ws.Range("E3").Value = "6,36"
$ws.Range("F3").Value = "11 de jun., 16:36 UTC ·"

# Row 4
$ws.Range("F4").Value = "11 de jun., 16:37 UTC ·"

# Row 5
$ws.Range("F5").Value = "11 de jun., 16:37 UTC ·"

# Row 6
$ws.Range("F6").Value = "11 de jun., 16:37 UTC ·"

# Row 7
$ws.Range("E7").Value = "6,75"
$ws.Range("F7").Value = "11 de jun., 16:36 UTC ·"

# Row 8
$ws.Range("E8").Value = "3,61"
$ws.Range("F8").Value = "11 de jun., 16:36 UTC ·"

# Row 9
$ws.Range("F9").Value = "11 de jun., 16:36 UTC ·"

# Row 10
$ws.Range("E10").Value = "4,05"
$ws.Range("F10").Value = "11 de jun., 16:36 UTC ·"

# Row 11
$ws.Range("F11").Value = "11 de jun., 16:36 UTC ·"

# Row 12
$ws.Range("F12").Value = "11 de jun., 16:37 UTC ·"

# Row 13
$ws.Range("F13").Value = "11 de jun., 16:37 UTC ·"

# Row 14 - no change

# Row 15
$ws.Range("F15").Value = "11 de jun., 16:36 UTC ·"

# Row 16
$ws.Range("F16").Value = "11 de jun., 16:36 UTC ·"

# Row 17
$ws.Range("F17").Value = "11 de jun., 16:37 UTC ·"

# Row 18
$ws.Range("F18").Value = "11 de jun., 16:37 UTC ·"

# Row 19
$ws.Range("E19").Value = "4,31"
$ws.Range("F19").Value = "11 de jun., 16:36 UTC ·"

# Row 20
$ws.Range("F20").Value = "11 de jun., 16:37 UTC ·"

# Row 21
$ws.Range("E21").Value = "4,15"
$ws.Range("F21").Value = "11 de jun., 16:37 UTC ·"

# Row 22
$ws.Range("F22").Value = "11 de jun., 16:37 UTC ·"

# Row 23
$ws.Range("F23").Value = "11 de jun., 16:36 UTC ·"

# Row 24
$ws.Range("F24").Value = "11 de jun., 16:37 UTC ·"

# Row 25
$ws.Range("F25").Value = "11 de jun., 16:37 UTC ·"
